$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "55.021.97"
$ws.Range("E2").Value = "  -1.70%  "
$ws.Range("D3").Value = "2.335.47"
$ws.Range("E3").Value = "  -5.03%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "'474.22"
$ws.Range("E5").Value = "  -2.57%  "
$ws.Range("D6").Value = "'144.08"
$ws.Range("E6").Value = "  -0.95%  "
$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "'0.610"
$ws.Range("E8").Value = "  +20.18%  "
$ws.Range("D9").Value = "2.332.03"
$ws.Range("E9").Value = "  -5.32%  "
$ws.Range("D10").Value = "'0.0953"
$ws.Range("E10").Value = "  -1.38%  "
$ws.Range("D11").Value = "'5.42"
$ws.Range("E11").Value = "  -6.34%  "
$ws.Range("D12").Value = "'0.325"
$ws.Range("E12").Value = "  -1.66%  "
$ws.Range("D13").Value = "'0.125"
$ws.Range("E13").Value = "  +1.19%  "
$ws.Range("D14").Value = "2.741.12"
$ws.Range("E14").Value = "  -4.93%  "
$ws.Range("D15").Value = "55.068.19"
$ws.Range("E15").Value = "  -1.63%  "
$ws.Range("D16").Value = "'19.88"
$ws.Range("E16").Value = "  -5.31%  "
$ws.Range("D17").Value = "'0.0000129"
$ws.Range("E17").Value = "  -5.11%  "
$ws.Range("D18").Value = "2.335.15"
$ws.Range("E18").Value = "  -5.32%  "
$ws.Range("D19").Value = "'4.54"
$ws.Range("E19").Value = "  +0.97%  "
$ws.Range("D20").Value = "'313.15"
$ws.Range("E20").Value = "  -0.91%  "
$ws.Range("D21").Value = "'9.52"
$ws.Range("E21").Value = "  -5.14%  "
$ws.Range("D22").Value = "'0.999"
$ws.Range("E22").Value = "  +0.18%  "
$ws.Range("D23").Value = "'5.61"
$ws.Range("E23").Value = "  -2.70%  "
$ws.Range("D24").Value = "'55.91"
$ws.Range("E24").Value = "  -4.05%  "
$ws.Range("E25").Value = "  -0.16%  "
$ws.Range("D26").Value = "'0.392"
$ws.Range("E26").Value = "  -4.74%  "
$ws.Range("B27").Value = "Kaspa"
$ws.Range("C27").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D27").Value = "'0.150"
$ws.Range("E27").Value = "  -5.67%  "
$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D28").Value = "2.427.78"
$ws.Range("E28").Value = "  -6.15%  "
$ws.Range("D29").Value = "'7.02"
$ws.Range("E29").Value = "  -8.32%  "
$ws.Range("E30").Value = "  +0.06%  "
$ws.Range("D31").Value = "0.0₃0736"
$ws.Range("E31").Value = "  -5.75%  "
$ws.Range("D32").Value = "'145.54"
$ws.Range("E32").Value = "  -1.53%  "
$ws.Range("D33").Value = "'18.03"
$ws.Range("E33").Value = "  -0.68%  "
$ws.Range("E34").Value = "  -2.30%  "
$ws.Range("D35").Value = "'5.06"
$ws.Range("E35").Value = "  -1.95%  "
$ws.Range("D36").Value = "'1.09"
$ws.Range("E36").Value = "  -4.37%  "
$ws.Range("D37").Value = "'3.58"
$ws.Range("E37").Value = "  -3.73%  "
$ws.Range("D38").Value = "'0.805"
$ws.Range("E38").Value = "  -6.64%  "
$ws.Range("E39").Value = "  +10.04%  "
$ws.Range("D40").Value = "'33.53"
$ws.Range("E40").Value = "  -0.67%  "
$ws.Range("E41").Value = "  +0.14%  "
$ws.Range("D42").Value = "'1.32"
$ws.Range("E42").Value = "  -0.41%  "
$ws.Range("D43").Value = "'3.36"
$ws.Range("E43").Value = "  -4.02%  "
$ws.Range("D44").Value = "'0.572"
$ws.Range("E44").Value = "  -4.94%  "
$ws.Range("B45").Value = "Hedera"
$ws.Range("C45").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D45").Value = "'0.0513"
$ws.Range("E45").Value = "  -7.20%  "
$ws.Range("B46").Value = "WhiteBITCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D46").Value = "'10.14"
$ws.Range("E46").Value = "  -0.41%  "
$ws.Range("D47").Value = "'248.37"
$ws.Range("E47").Value = "  -4.50%  "
$ws.Range("E48").Value = "  -3.38%  "
$ws.Range("D49").Value = "'4.34"
$ws.Range("E49").Value = "  -7.58%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'16.58"
$ws.Range("E50").Value = "  -4.92%  "
$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").Value = "1.779.73"
$ws.Range("E51").Value = "  -4.66%  "
